$d = $word.ActiveDocument

# Locate the paragraph holding the "m:'doc.html'.fromHTMLURI()" field
# (a Word field made of fldChar begin/end + instrText runs, split around
# a _GoBack bookmark) and rewrite it as plain literal text runs
# "{", "m", ":", "'", "doc.html", "'.fromHTMLURI()", "}" -- i.e. turn the
# field code into the literal M2Doc tag text, keeping the bookmark in
# place, per TokenIteratorFieldRewriterSplit.
$f = $d.Fields(1)
$p = $f.Code.Paragraphs(1)
$r = $p.Range

$q = "'"

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F">' + `
  '<w:r><w:t>{</w:t></w:r>' + `
  '<w:r><w:t>m</w:t></w:r>' + `
  '<w:r><w:t>:</w:t></w:r>' + `
  '<w:r><w:t>' + $q + '</w:t></w:r>' + `
  '<w:r><w:t>doc.html</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
  '<w:r><w:t>' + $q + '.fromHTMLURI()</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">}</w:t></w:r>' + `
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)
